$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update the existing "auto" prediction row with refreshed numbers ---
$ws.Cells.Item(2, 1).Value  = 42651.601539351854   # A2 Date
$ws.Cells.Item(2, 2).Value  = 24                    # B2 ScoreFinal
$ws.Cells.Item(2, 3).Value  = "Strong Buy"          # C2 Verdict
$ws.Cells.Item(2, 4).Value  = 32                    # D2 totalSentiment
$ws.Cells.Item(2, 5).Value  = 9592                  # E2 wordCount
$ws.Cells.Item(2, 6).Value  = 544                   # F2 sentenceCount
$ws.Cells.Item(2, 7).Value  = 52                    # G2 posWordPercentage
$ws.Cells.Item(2, 8).Value  = 47                    # H2 negWordPercentage
$ws.Cells.Item(2, 9).Value  = 90                    # I2 posPhrasePercentage
$ws.Cells.Item(2, 10).Value = 9                     # J2 negPhrasePercentage
$ws.Cells.Item(2, 11).Value = 63281                 # K2 ElapsedMs
$ws.Cells.Item(2, 12).Value = 100                   # L2 posWordCount
$ws.Cells.Item(2, 13).Value = 91                    # M2 negWordCount
$ws.Cells.Item(2, 14).Value = 50                    # N2 positivePhraseCount
$ws.Cells.Item(2, 15).Value = 5                     # O2 negativePhraseCount
$ws.Cells.Item(2, 16).Value = "Named"                # P2 Method
$ws.Cells.Item(2, 17).Value = 29.378539412357895    # Q2 RSI
$ws.Cells.Item(2, 18).Value = 0.84                  # R2 PEG
$ws.Cells.Item(2, 19).Value = -0.0136 # S2 200Moving%
$ws.Cells.Item(2, 20).Value = -0.03                 # T2 50Moving%
$ws.Cells.Item(2, 21).Value = 14.53                 # U2 PriceBook
$ws.Cells.Item(2, 22).Value = "N/A"                 # V2 Dividend
$ws.Cells.Item(2, 23).Value = 1                     # W2 Bollinger
$ws.Cells.Item(2, 24).Value = 0                     # X2 PriceChange (new)
$ws.Cells.Item(2, 25).Value = "Up"                  # Y2 UpDown (new)

# --- Row 3: brand new "named" row with mostly-empty sentiment numbers ---
$ws.Cells.Item(3, 1).Value  = 42651.601967592593    # A3 Date
$ws.Cells.Item(3, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(3, 2).Value  = 14                    # B3 ScoreFinal
$ws.Cells.Item(3, 3).Value  = "Buy"                 # C3 Verdict
$ws.Cells.Item(3, 4).Value  = 0                     # D3 totalSentiment
$ws.Cells.Item(3, 5).Value  = 0                     # E3 wordCount
$ws.Cells.Item(3, 6).Value  = 0                     # F3 sentenceCount
$ws.Cells.Item(3, 7).Value  = 0                     # G3 posWordPercentage
$ws.Cells.Item(3, 8).Value  = 0                     # H3 negWordPercentage
$ws.Cells.Item(3, 9).Value  = 0                     # I3 posPhrasePercentage
$ws.Cells.Item(3, 10).Value = 0                     # J3 negPhrasePercentage
$ws.Cells.Item(3, 11).Value = 4748                  # K3 ElapsedMs
$ws.Cells.Item(3, 12).Value = 0                     # L3 posWordCount
$ws.Cells.Item(3, 13).Value = 0                     # M3 negWordCount
$ws.Cells.Item(3, 14).Value = 0                     # N3 positivePhraseCount
$ws.Cells.Item(3, 15).Value = 0                     # O3 negativePhraseCount
$ws.Cells.Item(3, 16).Value = "Named"                # P3 Method
$ws.Cells.Item(3, 17).Value = 29.820796582770228    # Q3 RSI
$ws.Cells.Item(3, 18).Value = 0.84                  # R3 PEG
$ws.Cells.Item(3, 19).Value = -0.0136 # S3 200Moving%
$ws.Cells.Item(3, 19).NumberFormat = "0.00%"
$ws.Cells.Item(3, 20).Value = -0.03                 # T3 50Moving%
$ws.Cells.Item(3, 20).NumberFormat = "0.00%"
$ws.Cells.Item(3, 21).Value = 14.53                 # U3 PriceBook
$ws.Cells.Item(3, 22).Value = "N/A"                 # V3 Dividend
$ws.Cells.Item(3, 23).Value = 1                     # W3 Bollinger
$ws.Cells.Item(3, 24).Value = -0.39000000000000057  # X3 PriceChange
$ws.Cells.Item(3, 25).Value = "Down"                # Y3 UpDown

# --- Row 4: another new "named" row, repeating the refreshed auto numbers ---
$ws.Cells.Item(4, 1).Value  = 42651.682488425926    # A4 Date
$ws.Cells.Item(4, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(4, 2).Value  = 11                    # B4 ScoreFinal
$ws.Cells.Item(4, 3).Value  = "Buy"                 # C4 Verdict
$ws.Cells.Item(4, 4).Value  = 32                    # D4 totalSentiment
$ws.Cells.Item(4, 5).Value  = 9594                  # E4 wordCount
$ws.Cells.Item(4, 6).Value  = 544                   # F4 sentenceCount
$ws.Cells.Item(4, 7).Value  = 52                    # G4 posWordPercentage
$ws.Cells.Item(4, 8).Value  = 47                    # H4 negWordPercentage
$ws.Cells.Item(4, 9).Value  = 90                    # I4 posPhrasePercentage
$ws.Cells.Item(4, 10).Value = 9                     # J4 negPhrasePercentage
$ws.Cells.Item(4, 11).Value = 51453                 # K4 ElapsedMs
$ws.Cells.Item(4, 12).Value = 100                   # L4 posWordCount
$ws.Cells.Item(4, 13).Value = 91                    # M4 negWordCount
$ws.Cells.Item(4, 14).Value = 50                    # N4 positivePhraseCount
$ws.Cells.Item(4, 15).Value = 5                     # O4 negativePhraseCount
$ws.Cells.Item(4, 16).Value = "Named"                # P4 Method
$ws.Cells.Item(4, 17).Value = 29.009771469523784    # Q4 RSI
$ws.Cells.Item(4, 18).Value = 0.84                  # R4 PEG
$ws.Cells.Item(4, 19).Value = -0.0136 # S4 200Moving%
$ws.Cells.Item(4, 19).NumberFormat = "0.00%"
$ws.Cells.Item(4, 20).Value = -0.03                 # T4 50Moving%
$ws.Cells.Item(4, 20).NumberFormat = "0.00%"
$ws.Cells.Item(4, 21).Value = 14.53                 # U4 PriceBook
$ws.Cells.Item(4, 22).Value = "N/A"                 # V4 Dividend
$ws.Cells.Item(4, 23).Value = -2                    # W4 Bollinger

# --- Column C got wider now that "Strong Buy" lives in it ---
$ws.Columns.Item(3).ColumnWidth = 9.59
